# Apply the edit described by the diff:
# - Cell C6 changes from the numeric value 0.10854 to the text "0,10854 wie viel weniger?"
# - A new row 7 is added with: A7=5, B7="?", C7="?", D7="?", E7="q5"
# - Selection moves to G7

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 7 contents
$ws.Range("A7").Value = 5
$ws.Range("E7").Value = "q5"
$ws.Range("B7").Value = "?"
$ws.Range("C7").Value = "?"
$ws.Range("D7").Value = "?"

# Update C6 to hold text instead of the previous numeric value
$ws.Range("C6").Value = "0,10854 wie viel weniger?"

# Update the active selection to match the recorded state
$ws.Range("G7").Select()
